# Applies the values described by the diff for Sheets/Odin_Profits.xlsx
# (market-price / leve-profit recompute across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1555.4584
$ws.Range("I107").Value = 1592.2609
$ws.Range("K107").Value = 1592.2609
$ws.Range("M107").Value = 327.7391
$ws.Range("H132").Value = 343154.78
$ws.Range("I132").Value = 378672.06
$ws.Range("K132").Value = 1136016.18
$ws.Range("M132").Value = -1133486.18

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2711.7273
$ws.Range("I122").Value = 2132.9
$ws.Range("K122").Value = 6398.700000000001
$ws.Range("M122").Value = -3948.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 5547.5713
$ws.Range("I64").Value = 1674.5
$ws.Range("K64").Value = 1674.5
$ws.Range("M64").Value = -1449.5
$ws.Range("H67").Value = 5547.5713
$ws.Range("I67").Value = 1674.5
$ws.Range("K67").Value = 1674.5
$ws.Range("M67").Value = -894.5
$ws.Range("H82").Value = 18157.285
$ws.Range("J82").Value = 12345
$ws.Range("L82").Value = 12345
$ws.Range("N82").Value = -13111
$ws.Range("H85").Value = 18157.285
$ws.Range("J85").Value = 12345
$ws.Range("L85").Value = 12345
$ws.Range("N85").Value = -14997
$ws.Range("H107").Value = 2861519
$ws.Range("I107").Value = 3849450.8
$ws.Range("K107").Value = 3849450.8
$ws.Range("M107").Value = -3847530.8
$ws.Range("H134").Value = 2511523.8
$ws.Range("I134").Value = 2642261.8
$ws.Range("K134").Value = 7926785.399999999
$ws.Range("M134").Value = -7924250.399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4759.7393
$ws.Range("I31").Value = 1182.5555
$ws.Range("J31").Value = 7059.357
$ws.Range("K31").Value = 1182.5555
$ws.Range("L31").Value = 7059.357
$ws.Range("M31").Value = -887.5554999999999
$ws.Range("N31").Value = -7649.357
$ws.Range("H34").Value = 4759.7393
$ws.Range("I34").Value = 1182.5555
$ws.Range("J34").Value = 7059.357
$ws.Range("K34").Value = 1182.5555
$ws.Range("L34").Value = 7059.357
$ws.Range("M34").Value = -980.5554999999999
$ws.Range("N34").Value = -7463.357
$ws.Range("H94").Value = 38464560
$ws.Range("I94").Value = 76924776
$ws.Range("K94").Value = 76924776
$ws.Range("M94").Value = -76924325
$ws.Range("H105").Value = 62501276
$ws.Range("I105").Value = 62501276
$ws.Range("K105").Value = 62501276
$ws.Range("M105").Value = -62499529
$ws.Range("H107").Value = 1264.1538
$ws.Range("I107").Value = 858.1818
$ws.Range("J107").Value = 3497
$ws.Range("K107").Value = 858.1818
$ws.Range("L107").Value = 3497
$ws.Range("M107").Value = 1061.8182
$ws.Range("N107").Value = -7337
$ws.Range("H132").Value = 9993.585999999999
$ws.Range("J132").Value = 15161
$ws.Range("L132").Value = 45483
$ws.Range("N132").Value = -50543
$ws.Range("H134").Value = 71435270
$ws.Range("I134").Value = 83338260
$ws.Range("J134").Value = 17375
$ws.Range("K134").Value = 250014780
$ws.Range("L134").Value = 52125
$ws.Range("M134").Value = -250012245
$ws.Range("N134").Value = -57195

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 999.5
$ws.Range("I50").Value = 999.5
$ws.Range("K50").Value = 2998.5
$ws.Range("M50").Value = -2517.5
$ws.Range("H53").Value = 999.5
$ws.Range("I53").Value = 999.5
$ws.Range("K53").Value = 2998.5
$ws.Range("M53").Value = -2517.5
$ws.Range("H108").Value = 6625
$ws.Range("I108").Value = 6625
$ws.Range("K108").Value = 19875
$ws.Range("M108").Value = -16995
$ws.Range("H132").Value = 2132.1667
$ws.Range("I132").Value = 937.5
$ws.Range("J132").Value = 3087.9
$ws.Range("K132").Value = 8437.5
$ws.Range("L132").Value = 27791.1
$ws.Range("M132").Value = -5907.5
$ws.Range("N132").Value = -32851.10000000001
$ws.Range("H139").Value = 3509.5386
$ws.Range("I139").Value = 2161.5715
$ws.Range("J139").Value = 5082.1665
$ws.Range("K139").Value = 6484.7145
$ws.Range("L139").Value = 15246.4995
$ws.Range("M139").Value = -1344.7145
$ws.Range("N139").Value = -25526.4995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H112").Value = 75000
$ws.Range("J112").Value = 75000
$ws.Range("L112").Value = 75000
$ws.Range("N112").Value = -77216
$ws.Range("H113").Value = 5015.9287
$ws.Range("I113").Value = 1780.7222
$ws.Range("K113").Value = 1780.7222
$ws.Range("M113").Value = 389.2778000000001
$ws.Range("H126").Value = 12826738
$ws.Range("I126").Value = 20003252
$ws.Range("J126").Value = 11534.929
$ws.Range("K126").Value = 60009756
$ws.Range("L126").Value = 34604.787
$ws.Range("M126").Value = -60007286
$ws.Range("N126").Value = -39544.787
$ws.Range("H132").Value = 27031286
$ws.Range("I132").Value = 40003884
$ws.Range("J132").Value = 5046.6665
$ws.Range("K132").Value = 120011652
$ws.Range("L132").Value = 15139.9995
$ws.Range("M132").Value = -120009122
$ws.Range("N132").Value = -20199.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6844.276
$ws.Range("I7").Value = 6369.2383
$ws.Range("J7").Value = 8091.25
$ws.Range("K7").Value = 6369.2383
$ws.Range("L7").Value = 8091.25
$ws.Range("M7").Value = -6257.2383
$ws.Range("N7").Value = -8315.25
$ws.Range("H22").Value = 915.53845
$ws.Range("I22").Value = 769.1667
$ws.Range("J22").Value = 1041
$ws.Range("K22").Value = 769.1667
$ws.Range("L22").Value = 1041
$ws.Range("M22").Value = -474.1667
$ws.Range("N22").Value = -1631
$ws.Range("H27").Value = 915.53845
$ws.Range("I27").Value = 769.1667
$ws.Range("J27").Value = 1041
$ws.Range("K27").Value = 769.1667
$ws.Range("L27").Value = 1041
$ws.Range("M27").Value = -662.1667
$ws.Range("N27").Value = -1255
$ws.Range("H55").Value = 4776.9443
$ws.Range("I55").Value = 2625
$ws.Range("J55").Value = 6498.5
$ws.Range("K55").Value = 2625
$ws.Range("L55").Value = 6498.5
$ws.Range("M55").Value = -2452
$ws.Range("N55").Value = -6844.5
$ws.Range("H122").Value = 5887.3
$ws.Range("I122").Value = 5333
$ws.Range("K122").Value = 15999
$ws.Range("M122").Value = -13549
$ws.Range("H126").Value = 6844.276
$ws.Range("I126").Value = 6369.2383
$ws.Range("J126").Value = 8091.25
$ws.Range("K126").Value = 19107.7149
$ws.Range("L126").Value = 24273.75
$ws.Range("M126").Value = -16637.7149
$ws.Range("N126").Value = -29213.75
$ws.Range("H132").Value = 10103.577
$ws.Range("I132").Value = 9537.111000000001
$ws.Range("K132").Value = 28611.333
$ws.Range("M132").Value = -26081.333
$ws.Range("H136").Value = 15159910
$ws.Range("I136").Value = 27784182
$ws.Range("J136").Value = 10782.733
$ws.Range("K136").Value = 83352546
$ws.Range("L136").Value = 32348.199
$ws.Range("M136").Value = -83349996
$ws.Range("N136").Value = -37448.199

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 12518000
$ws.Range("I32").Value = 12518000
$ws.Range("K32").Value = 12518000
$ws.Range("M32").Value = -12517683
$ws.Range("H34").Value = 37990
$ws.Range("I34").Value = 37990
$ws.Range("K34").Value = 37990
$ws.Range("M34").Value = -37787
$ws.Range("H107").Value = 5883173
$ws.Range("I107").Value = 9091500
$ws.Range("K107").Value = 27274500
$ws.Range("M107").Value = -27272580
$ws.Range("H122").Value = 12394
$ws.Range("I122").Value = 7019.1
$ws.Range("K122").Value = 21057.3
$ws.Range("M122").Value = -18607.3
$ws.Range("H132").Value = 7577.44
$ws.Range("I132").Value = 6485.727
$ws.Range("K132").Value = 19457.181
$ws.Range("M132").Value = -16927.181
